{"js": "// Each (oldText, newText) pair corresponds to one equation cell in the\n// worksheet table (\"before -> after\" from the target diff). Every\n// \"oldText\" is unique in the document, and matchWholeWord keeps a short\n// equation like \"0+18=\" from matching inside a longer one such as\n// \"40+18=\" that happens to contain it as a substring.\nconst pairs = [\n  [\"3+44=\", \"70-7=\"],\n  [\"86-54=\", \"17+55=\"],\n  [\"55+8=\", \"67-42=\"],\n  [\"37+56=\", \"25+58=\"],\n  [\"45-38=\", \"11+45=\"],\n  [\"56-32=\", \"34+9=\"],\n  [\"17+3=\", \"58-36=\"],\n  [\"30+44=\", \"96-62=\"],\n  [\"19-10=\", \"46-13=\"],\n  [\"43+28=\", \"39+3=\"],\n  [\"30+20=\", \"99-65=\"],\n  [\"94-51=\", \"65+21=\"],\n  [\"62+12=\", \"31-17=\"],\n  [\"64+27=\", \"19+38=\"],\n  [\"74-69=\", \"5+48=\"],\n  [\"66-2=\", \"10+2=\"],\n  [\"91-52=\", \"57-50=\"],\n  [\"4+51=\", \"84-76=\"],\n  [\"5+19=\", \"10+89=\"],\n  [\"70-51=\", \"85-19=\"],\n  [\"97-69=\", \"74+15=\"],\n  [\"68+14=\", \"62+18=\"],\n  [\"15+60=\", \"21+56=\"],\n  [\"51-39=\", \"60+11=\"],\n  [\"21+11=\", \"95-85=\"],\n  [\"31+18=\", \"36+23=\"],\n  [\"71-30=\", \"54+4=\"],\n  [\"64-41=\", \"19+48=\"],\n  [\"95-29=\", \"9+49=\"],\n  [\"61+7=\", \"82-40=\"],\n  [\"41-25=\", \"9+2=\"],\n  [\"3+43=\", \"11+17=\"],\n  [\"50+38=\", \"36+47=\"],\n  [\"30+30=\", \"28+47=\"],\n  [\"3+83=\", \"77-33=\"],\n  [\"16+37=\", \"47+1=\"],\n  [\"74-29=\", \"15+21=\"],\n  [\"84+2=\", \"84-21=\"],\n  [\"10+28=\", \"24-12=\"],\n  [\"86-53=\", \"99-42=\"],\n  [\"22-10=\", \"4+4=\"],\n  [\"67-57=\", \"75-62=\"],\n  [\"50-40=\", \"78-18=\"],\n  [\"34+39=\", \"94-71=\"],\n  [\"96-51=\", \"38+11=\"],\n  [\"69-47=\", \"97+1=\"],\n  [\"31-4=\", \"92-86=\"],\n  [\"62-17=\", \"10+16=\"],\n  [\"34+36=\", \"78-25=\"],\n  [\"18+16=\", \"20+26=\"],\n  [\"34-1=\", \"87-24=\"],\n  [\"27-8=\", \"33+58=\"],\n  [\"24-14=\", \"53+12=\"],\n  [\"93-42=\", \"38+9=\"],\n  [\"24+67=\", \"84-76=\"],\n  [\"3+85=\", \"54+16=\"],\n  [\"7+84=\", \"61-13=\"],\n  [\"50-18=\", \"34-14=\"],\n  [\"23+70=\", \"20+62=\"],\n  [\"25+10=\", \"5+74=\"],\n  [\"46+41=\", \"98-47=\"],\n  [\"38+41=\", \"89-32=\"],\n  [\"79-44=\", \"18+30=\"],\n  [\"95-34=\", \"75-32=\"],\n  [\"56-50=\", \"50-47=\"],\n  [\"44-3=\", \"55+31=\"],\n  [\"63+30=\", \"60+6=\"],\n  [\"47+7=\", \"81-65=\"],\n  [\"84-70=\", \"92-90=\"],\n  [\"85-47=\", \"98-92=\"],\n  [\"35+9=\", \"30-18=\"],\n  [\"84-4=\", \"22+33=\"],\n  [\"0+18=\", \"10+80=\"],\n  [\"40+2=\", \"78-15=\"],\n  [\"67-61=\", \"51-0=\"],\n  [\"69-33=\", \"9+2=\"],\n  [\"55+12=\", \"87-31=\"],\n  [\"64+6=\", \"33+9=\"],\n  [\"40+18=\", \"9+45=\"],\n  [\"25+2=\", \"4+21=\"],\n  [\"56-13=\", \"50-41=\"],\n  [\"83-73=\", \"39-27=\"],\n  [\"76-40=\", \"56-2=\"],\n  [\"48+0=\", \"45+9=\"],\n  [\"56+13=\", \"48-26=\"],\n  [\"48+39=\", \"4+48=\"],\n  [\"44-29=\", \"21-18=\"],\n  [\"79-67=\", \"13+78=\"],\n  [\"22-4=\", \"44-28=\"],\n  [\"77-70=\", \"75-67=\"],\n  [\"16-0=\", \"70-8=\"],\n  [\"38-21=\", \"59-23=\"],\n  [\"44-42=\", \"58-57=\"],\n  [\"76+4=\", \"30+18=\"],\n  [\"95-11=\", \"12+39=\"],\n  [\"63+35=\", \"75-20=\"],\n  [\"77+16=\", \"50-27=\"],\n  [\"74-54=\", \"10+71=\"],\n  [\"22-7=\", \"78-42=\"],\n  [\"50-29=\", \"91-45=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(newText, \"Replace\");\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Each (old, new) pair below corresponds to one equation cell in the\n# worksheet table. Every \"old\" string is unique in the document. We use\n# MatchWholeWord so e.g. replacing \"0+18=\" cannot also clobber the\n# \"40+18=\" cell that happens to contain \"0+18=\" as a substring.\n$pairs = @(\n    ,@('3+44=', '70-7=')\n    ,@('86-54=', '17+55=')\n    ,@('55+8=', '67-42=')\n    ,@('37+56=', '25+58=')\n    ,@('45-38=', '11+45=')\n    ,@('56-32=', '34+9=')\n    ,@('17+3=', '58-36=')\n    ,@('30+44=', '96-62=')\n    ,@('19-10=', '46-13=')\n    ,@('43+28=', '39+3=')\n    ,@('30+20=', '99-65=')\n    ,@('94-51=', '65+21=')\n    ,@('62+12=', '31-17=')\n    ,@('64+27=', '19+38=')\n    ,@('74-69=', '5+48=')\n    ,@('66-2=', '10+2=')\n    ,@('91-52=', '57-50=')\n    ,@('4+51=', '84-76=')\n    ,@('5+19=', '10+89=')\n    ,@('70-51=', '85-19=')\n    ,@('97-69=', '74+15=')\n    ,@('68+14=', '62+18=')\n    ,@('15+60=', '21+56=')\n    ,@('51-39=', '60+11=')\n    ,@('21+11=', '95-85=')\n    ,@('31+18=', '36+23=')\n    ,@('71-30=', '54+4=')\n    ,@('64-41=', '19+48=')\n    ,@('95-29=', '9+49=')\n    ,@('61+7=', '82-40=')\n    ,@('41-25=', '9+2=')\n    ,@('3+43=', '11+17=')\n    ,@('50+38=', '36+47=')\n    ,@('30+30=', '28+47=')\n    ,@('3+83=', '77-33=')\n    ,@('16+37=', '47+1=')\n    ,@('74-29=', '15+21=')\n    ,@('84+2=', '84-21=')\n    ,@('10+28=', '24-12=')\n    ,@('86-53=', '99-42=')\n    ,@('22-10=', '4+4=')\n    ,@('67-57=', '75-62=')\n    ,@('50-40=', '78-18=')\n    ,@('34+39=', '94-71=')\n    ,@('96-51=', '38+11=')\n    ,@('69-47=', '97+1=')\n    ,@('31-4=', '92-86=')\n    ,@('62-17=', '10+16=')\n    ,@('34+36=', '78-25=')\n    ,@('18+16=', '20+26=')\n    ,@('34-1=', '87-24=')\n    ,@('27-8=', '33+58=')\n    ,@('24-14=', '53+12=')\n    ,@('93-42=', '38+9=')\n    ,@('24+67=', '84-76=')\n    ,@('3+85=', '54+16=')\n    ,@('7+84=', '61-13=')\n    ,@('50-18=', '34-14=')\n    ,@('23+70=', '20+62=')\n    ,@('25+10=', '5+74=')\n    ,@('46+41=', '98-47=')\n    ,@('38+41=', '89-32=')\n    ,@('79-44=', '18+30=')\n    ,@('95-34=', '75-32=')\n    ,@('56-50=', '50-47=')\n    ,@('44-3=', '55+31=')\n    ,@('63+30=', '60+6=')\n    ,@('47+7=', '81-65=')\n    ,@('84-70=', '92-90=')\n    ,@('85-47=', '98-92=')\n    ,@('35+9=', '30-18=')\n    ,@('84-4=', '22+33=')\n    ,@('0+18=', '10+80=')\n    ,@('40+2=', '78-15=')\n    ,@('67-61=', '51-0=')\n    ,@('69-33=', '9+2=')\n    ,@('55+12=', '87-31=')\n    ,@('64+6=', '33+9=')\n    ,@('40+18=', '9+45=')\n    ,@('25+2=', '4+21=')\n    ,@('56-13=', '50-41=')\n    ,@('83-73=', '39-27=')\n    ,@('76-40=', '56-2=')\n    ,@('48+0=', '45+9=')\n    ,@('56+13=', '48-26=')\n    ,@('48+39=', '4+48=')\n    ,@('44-29=', '21-18=')\n    ,@('79-67=', '13+78=')\n    ,@('22-4=', '44-28=')\n    ,@('77-70=', '75-67=')\n    ,@('16-0=', '70-8=')\n    ,@('38-21=', '59-23=')\n    ,@('44-42=', '58-57=')\n    ,@('76+4=', '30+18=')\n    ,@('95-11=', '12+39=')\n    ,@('63+35=', '75-20=')\n    ,@('77+16=', '50-27=')\n    ,@('74-54=', '10+71=')\n    ,@('22-7=', '78-42=')\n    ,@('50-29=', '91-45=')\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #         MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n    #         Format, ReplaceWith, Replace:=wdReplaceAll)\n    $null = $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
